# Trading update: 2026-02-18 11:00:15
#
# - Trade #9 (row 10 on "All Trades", row 2 on "MarketMaking") finished
#   filling / was re-priced: time, side and entry price changed.
# - Trade #10 was opened and appended as a new row to both the
#   "All Trades" log and the per-strategy "MarketMaking" sheet.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, [int]$row, [int]$col, [string]$text) {
    # Leading apostrophe forces Excel to store the value as literal text
    # instead of auto-converting date/time/number-looking strings (or
    # collapsing an empty string into a truly blank cell). Resetting the
    # style back to "Normal" afterwards drops the quote-prefix formatting
    # flag that the apostrophe trick leaves behind, so the cell ends up
    # with plain, unstyled text - same as a freshly authored cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

function Update-ExistingTrade($ws, [int]$row) {
    Set-TextValue $ws $row 3 "10:58:46"            # C: Time
    $ws.Cells.Item($row, 5).Value = "UP"           # E: Side
    $ws.Cells.Item($row, 6).Value = 0.58           # F: Entry Price
}

function Add-NewTrade($ws, [int]$row) {
    $ws.Cells.Item($row, 1).Value = 10                          # A: Trade #
    Set-TextValue $ws $row 2 "2026-02-18"                       # B: Date
    Set-TextValue $ws $row 3 "10:58:52"                         # C: Time
    $ws.Cells.Item($row, 4).Value = "MarketMaking"              # D: Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"                      # E: Side
    $ws.Cells.Item($row, 6).Value = 0.45                        # F: Entry Price
    Set-TextValue $ws $row 7 ""                                 # G: Exit Price
    $ws.Cells.Item($row, 8).Value = "OPEN"                      # H: Status
    $ws.Cells.Item($row, 9).Value = 0                           # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0                          # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100                        # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                          # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                          # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                        # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 202 bps"  # O: Entry Reason
    Set-TextValue $ws $row 16 ""                                # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0                          # Q: Duration (min)
}

# --- "All Trades" sheet: update trade #9 (row 10) and append trade #10 (row 11) ---
$wsAll = $wb.Worksheets.Item("All Trades")
Update-ExistingTrade $wsAll 10
Add-NewTrade $wsAll 11

# --- "MarketMaking" strategy sheet: same update mirrored onto its own log ---
$wsStrategy = $wb.Worksheets.Item("MarketMaking")
Update-ExistingTrade $wsStrategy 2
Add-NewTrade $wsStrategy 3
